# Update the forecast-error table: refresh existing Q0..Q8 rows (2-10)
# with new values, and append a new Q9 row (row 11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..11 (columns B..G), corresponding to labels Q0..Q9
$data = @(
    @(0.4985553813290332, 2.131733426832194,  17.05044517982819, 4.129218470828128, 4.139797756631062, 51),
    @(0.3670740509330708, 2.310008463032842,  18.22424447820178, 4.268986352543397, 4.296355968907215, 50),
    @(0.6013248752497656, 2.260079145283003,  17.85224903368007, 4.225192188963725, 4.225523177608233, 49),
    @(0.3757086490286848, 2.257680262786419,  19.07933980741992, 4.36799036256033,  4.397854349388155, 48),
    @(0.6368657846289102, 2.180689379475908,  18.26526036982974, 4.27378759063079,  4.271757936889329, 47),
    @(0.2352128511284962, 1.884428425708629,  15.78161112617242, 3.972607597809331, 4.009458699620951, 46),
    @(0.4131478667422561, 1.866528751697462,  15.52343479139494, 3.939979034385201, 3.9625332330756,   45),
    @(0.2445497104634806, 1.984813716461438,  16.20214174503358, 4.02518841112234,  4.064202304816381, 44),
    @(0.4887008263018376, 1.804638728127334,  15.84511544644097, 3.980592348688945, 3.997232193593323, 43),
    @(0.2035882577814556, 1.941069045450287,  16.47627633247782, 4.059097970298058, 4.103130222273781, 42)
)

$labels = @("Q0","Q1","Q2","Q3","Q4","Q5","Q6","Q7","Q8","Q9")

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $labels[$i]
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
    $ws.Cells.Item($row, 6).Value = $data[$i][4]
    $ws.Cells.Item($row, 7).Value = $data[$i][5]
}

# Copy the styling from the A10 (existing labeled cell) to the new A11 label cell
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
